$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 95
$ws.Range("A95").Value = 14552579
$ws.Range("B95").Value = "'2025-09-03"
$ws.Range("C95").Value = "Tom Gentzsch"
$ws.Range("D95").Value = "Francesco Forti"
$ws.Range("E95").Value = "Gana Francesco Forti"
$ws.Range("F95").Value = 2.25
$ws.Range("G95").Value = "'"
$ws.Range("H95").Value = "'"

# Row 96
$ws.Range("A96").Value = 14551799
$ws.Range("B96").Value = "'2025-09-03"
$ws.Range("C96").Value = "David Jorda Sanchis"
$ws.Range("D96").Value = "Carlos Lopez Montagud"
$ws.Range("E96").Value = "Gana Carlos Lopez Montagud"
$ws.Range("F96").Value = 2.25
$ws.Range("G96").Value = "'"
$ws.Range("H96").Value = "'"
